$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.953.45'
$ws.Range("E2").Value = '  -1.46%  '
$ws.Range("D3").Value = '3.383.37'
$ws.Range("E3").Value = '  -0.70%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '573.79'
$ws.Range("E5").Value = '  -0.62%  '
$ws.Range("D6").Value = '136.82'
$ws.Range("E6").Value = '  -0.56%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.380.49'
$ws.Range("E8").Value = '  -0.77%  '
$ws.Range("E9").Value = '  -1.35%  '
$ws.Range("E10").Value = '  +1.71%  '
$ws.Range("E11").Value = '  -2.45%  '
$ws.Range("D12").Value = '0.388'
$ws.Range("E12").Value = '  -1.27%  '
$ws.Range("D13").Value = '3.962.08'
$ws.Range("E13").Value = '  -0.82%  '
$ws.Range("E14").Value = '  +0.74%  '
$ws.Range("D15").Value = '26.28'
$ws.Range("E15").Value = '  +3.36%  '
$ws.Range("D16").Value = '0.0000173'
$ws.Range("E16").Value = '  -2.89%  '
$ws.Range("D17").Value = '3.384.97'
$ws.Range("E17").Value = '  -0.75%  '
$ws.Range("D18").Value = '61.081.98'
$ws.Range("E18").Value = '  -1.25%  '
$ws.Range("D19").Value = '14.03'
$ws.Range("E19").Value = '  -1.40%  '
$ws.Range("D20").Value = '5.84'
$ws.Range("E20").Value = '  -0.76%  '
$ws.Range("D21").Value = '9.47'
$ws.Range("E21").Value = '  -0.35%  '
$ws.Range("D22").Value = '377.71'
$ws.Range("E22").Value = '  -3.07%  '
$ws.Range("D23").Value = '0.555'
$ws.Range("E23").Value = '  -2.88%  '
$ws.Range("D24").Value = '3.527.23'
$ws.Range("E24").Value = '  -0.50%  '
$ws.Range("E25").Value = '  +0.17%  '
$ws.Range("D26").Value = '0.0000125'
$ws.Range("D27").Value = '71.24'
$ws.Range("E27").Value = '  -0.42%  '
$ws.Range("D28").Value = '1.76'
$ws.Range("E28").Value = '  +10.50%  '
$ws.Range("D29").Value = '0.173'
$ws.Range("E29").Value = '  +8.29%  '
$ws.Range("D30").Value = '7.54'
$ws.Range("E30").Value = '  -2.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  -0.28%  '
$ws.Range("D32").Value = '8.17'
$ws.Range("E32").Value = '  -1.48%  '
$ws.Range("D33").Value = '2.16'
$ws.Range("E33").Value = '  -0.48%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").Value = '23.78'
$ws.Range("E35").Value = '  +0.94%  '
$ws.Range("E36").Value = '  -4.13%  '
$ws.Range("D37").Value = '6.87'
$ws.Range("E37").Value = '  -1.65%  '
$ws.Range("E38").Value = '  -0.66%  '
$ws.Range("D39").Value = '164.77'
$ws.Range("E39").Value = '  +0.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0760'
$ws.Range("E40").Value = '  -3.78%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("E42").Value = '  -2.60%  '
$ws.Range("E43").Value = '  -1.89%  '
$ws.Range("E44").Value = '  -1.02%  '
$ws.Range("D45").Value = '41.62'
$ws.Range("E45").Value = '  -0.10%  '
$ws.Range("E46").Value = '  -2.59%  '
$ws.Range("D47").Value = '24.09'
$ws.Range("E47").Value = '  -3.01%  '
$ws.Range("D48").Value = '2.478.87'
$ws.Range("E48").Value = '  +4.09%  '
$ws.Range("D49").Value = '23.23'
$ws.Range("E49").Value = '  -0.08%  '
$ws.Range("D50").Value = '6.81'
$ws.Range("E50").Value = '  -2.40%  '
$ws.Range("E51").Value = '  +3.97%  '
